$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "'2021-12-12"
$ws.Range("B2").Style = "Normal"
